$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table cell margin: left margin 118 -> 123 dxa (5.9pt -> 6.15pt)
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.LeftPadding = 6.15

# ---------------------------------------------------------------------------
# 2) Merge the date text + "г." runs into a single run with the new
#    placeholder name (humanized_created_at_with_month_as_word ->
#    humanized_created_at_with_quotes_and_month_as_word), dropping the
#    separate trailing run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "{issue.humanized_created_at_with_month_as_word} г.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{issue.humanized_created_at_with_quotes_and_month_as_word} г.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Add new character styles ListLabel30 .. ListLabel40 to styles.xml
#    (note: this interpreter only reliably binds POSITIONAL function
#    parameters, so avoid named -Param style calls here)
# ---------------------------------------------------------------------------
function Add-ListLabelStyle($Id, $DisplayName, $Ascii, $Cs, $SizePt, $Bold) {
    $style = $d.Styles.Add($Id, 2)
    $style.NameLocal = $DisplayName
    $style.QuickStyle = $true
    if ($Ascii -ne $null) {
        $style.Font.Name = $Ascii
    }
    if ($Cs -ne $null) {
        $style.Font.NameBi = $Cs
    }
    if ($SizePt -ne $null) {
        $style.Font.Size = $SizePt
    }
    if ($Bold -ne $null) {
        $style.Font.Bold = $Bold
    }
    return $style
}

Add-ListLabelStyle "ListLabel30" "ListLabel 30" "Times New Roman" $null 10.5 $true   | Out-Null
Add-ListLabelStyle "ListLabel31" "ListLabel 31" $null             $null $null $false | Out-Null
Add-ListLabelStyle "ListLabel32" "ListLabel 32" "Times New Roman" "Symbol"      10.5 $null | Out-Null
Add-ListLabelStyle "ListLabel33" "ListLabel 33" $null             "Courier New" $null $null | Out-Null
Add-ListLabelStyle "ListLabel34" "ListLabel 34" $null             "Wingdings"   $null $null | Out-Null
Add-ListLabelStyle "ListLabel35" "ListLabel 35" $null             "Symbol"      $null $null | Out-Null
Add-ListLabelStyle "ListLabel36" "ListLabel 36" $null             "Courier New" $null $null | Out-Null
Add-ListLabelStyle "ListLabel37" "ListLabel 37" $null             "Wingdings"   $null $null | Out-Null
Add-ListLabelStyle "ListLabel38" "ListLabel 38" $null             "Symbol"      $null $null | Out-Null
Add-ListLabelStyle "ListLabel39" "ListLabel 39" $null             "Courier New" $null $null | Out-Null
Add-ListLabelStyle "ListLabel40" "ListLabel 40" $null             "Wingdings"   $null $null | Out-Null

Write-Host "Done applying edits"
